$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells GH1 / GI1 (shared strings "Act_Impr" / "Act_Impr_2") ---
$ws.Range("GH1").Value = "Act_Impr"
$ws.Range("GI1").Value = "Act_Impr_2"

# --- New data values in columns GH (190) and GI (191) for rows 2-32 ---
$ws.Range("GH2").Value = 17
$ws.Range("GI2").Value = 22
$ws.Range("GH3").Value = 173
$ws.Range("GI3").Value = 56
$ws.Range("GH4").Value = 559
$ws.Range("GI4").Value = 234
$ws.Range("GH5").Value = 10
$ws.Range("GI5").Value = 6
$ws.Range("GH6").Value = 46
$ws.Range("GI6").Value = 29
$ws.Range("GH7").Value = 21
$ws.Range("GI7").Value = 9
$ws.Range("GH8").Value = 13
$ws.Range("GI8").Value = 7
$ws.Range("GH9").Value = 18
$ws.Range("GI9").Value = 5
$ws.Range("GH10").Value = 9
$ws.Range("GI10").Value = 7
$ws.Range("GH11").Value = 12
$ws.Range("GI11").Value = 3
$ws.Range("GH12").Value = 21
$ws.Range("GI12").Value = 8
$ws.Range("GH13").Value = 384
$ws.Range("GI13").Value = 158
$ws.Range("GH14").Value = 1
$ws.Range("GI14").Value = 2
$ws.Range("GH15").Value = 19
$ws.Range("GI15").Value = 7
$ws.Range("GH16").Value = 19
$ws.Range("GI16").Value = 4
$ws.Range("GH17").Value = 0
$ws.Range("GI17").Value = 0
$ws.Range("GH18").Value = 9
$ws.Range("GI18").Value = 3
$ws.Range("GH19").Value = 3
$ws.Range("GI19").Value = 1
$ws.Range("GH20").Value = 17
$ws.Range("GI20").Value = 6
$ws.Range("GH21").Value = 4
$ws.Range("GI21").Value = 3
$ws.Range("GH22").Value = 1
$ws.Range("GI22").Value = 3
$ws.Range("GH23").Value = 0
$ws.Range("GI23").Value = 1
$ws.Range("GH24").Value = 44
$ws.Range("GI24").Value = 13
$ws.Range("GH25").Value = 0
$ws.Range("GI25").Value = 2
$ws.Range("GH26").Value = 5
$ws.Range("GI26").Value = 0
$ws.Range("GH27").Value = 20
$ws.Range("GI27").Value = 13
$ws.Range("GH28").Value = 82
$ws.Range("GI28").Value = 31
$ws.Range("GH29").Value = 2
$ws.Range("GI29").Value = 0
$ws.Range("GH30").Value = 19
$ws.Range("GI30").Value = 11
$ws.Range("GH31").Value = 1
$ws.Range("GI31").Value = 0
$ws.Range("GH32").Value = 0
$ws.Range("GI32").Value = 0

# --- Highlight BG2 (formula cell) with the yellow "changed" fill ---
$ws.Range("BG2").Interior.Color = 65535

# --- Update the active selection to the new last column (GI1), matching the
#     view state left behind by the author's editing session ---
$ws.Activate()
$ws.Range("GI1").Select()
